# feat: add 2022-Q4 data
#
# 1. "总计" (summary) sheet gets a new top data row for 2022-Q4, with the
#    previously-existing rows shifted down by one and re-indexed.
# 2. A brand-new "2022-Q4" worksheet is inserted right after "总计" (i.e.
#    before the former first quarter sheet), carrying the fund holdings for
#    that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" summary sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$summaryRows = @(
    @(0, "2022-Q4", 7, 1.1),
    @(1, "2022-Q3", 2, 0.47),
    @(2, "2022-Q2", 2, 1.05),
    @(3, "2021-Q4", 1, 0.01),
    @(4, "2021-Q3", 1, 0.01)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
    $row = $i + 2
    $data = $summaryRows[$i]
    $total.Range("A$row").Value2 = $data[0]
    $total.Range("B$row").Value2 = $data[1]
    $total.Range("C$row").Value2 = $data[2]
    $total.Range("D$row").Value2 = $data[3]
}

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q4" sheet before the current first-quarter sheet
#    (copying it so the header/body formatting matches the other tabs).
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($q3)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

$fundRows = @(
    @(0, "007216", "浙商中华预期高股息C", "5.13", "91.62", "8.70", "0.4463", 1),
    @(1, "007178", "浙商中华预期高股息A", "3.53", "91.62", "8.70", "0.3071", 1),
    @(2, "008704", "广发高股息优享混合A", "2.52", "75.31", "5.53", "0.1394", 6),
    @(3, "013334", "鹏华价值远航6个月持有期混合A", "1.53", "93.83", "8.39", "0.1284", 5),
    @(4, "008705", "广发高股息优享混合C", "0.73", "75.31", "5.53", "0.0404", 6),
    @(5, "004292", "鹏华沪深港互联网股票", "0.76", "93.59", "3.77", "0.0287", 6),
    @(6, "013335", "鹏华价值远航6个月持有期混合C", "0.10", "93.83", "8.39", "0.0084", 5)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $row = $i + 2
    $data = $fundRows[$i]

    if ($row -gt 3) {
        # Rows 4-8 don't exist yet on the copied sheet - stamp the same
        # formatting used by the existing "A" column cells before writing.
        $newSheet.Range("A3").Copy()
        $newSheet.Range("A$row").PasteSpecial(-4122)
    }

    $newSheet.Range("A$row").Value2 = $data[0]
    $newSheet.Range("B$row").Value2 = $data[1]
    $newSheet.Range("C$row").Value2 = $data[2]
    $newSheet.Range("D$row").Value2 = $data[3]
    $newSheet.Range("E$row").Value2 = $data[4]
    $newSheet.Range("F$row").Value2 = $data[5]
    $newSheet.Range("G$row").Value2 = $data[6]
    $newSheet.Range("H$row").Value2 = $data[7]
}
